# PMD - UseStringBufferLength specified
# B31 currently reads "later" for the UseStringBufferLength rule; it has now
# been specified, so mark it "ok". This also ripples through the COUNTIF
# summary formulas in D2/D4 (and their dependent percentages in E2/E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B31").Value = "ok"

# Keep the active selection consistent with where the edit was made.
$ws.Range("B27").Select()

$wb.Application.Calculate()
